$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2 through 112 in column C (Fitness) had values of 7310 or 7295;
# update them all to 7293, matching the already-7293 values in rows 113+.
$ws.Range("C2:C112").Value = 7293
